$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 27000
$ws.Range("I21").Value = 24000
$ws.Range("J21").Value = 29000
$ws.Range("K21").Value = 24000
$ws.Range("L21").Value = 29000
$ws.Range("M21").Value = -23532
$ws.Range("N21").Value = -29936
$ws.Range("H23").Value = 27000
$ws.Range("I23").Value = 24000
$ws.Range("J23").Value = 29000
$ws.Range("K23").Value = 24000
$ws.Range("L23").Value = 29000
$ws.Range("M23").Value = -23766
$ws.Range("N23").Value = -29468
$ws.Range("H42").Value = 5565.5
$ws.Range("I42").Value = 5565.5
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 16696.5
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -16466.5
$ws.Range("N42").ClearContents()
$ws.Range("H62").Value = 5501.8237
$ws.Range("I62").Value = 5569.4443
$ws.Range("J62").Value = 5425.75
$ws.Range("K62").Value = 5569.4443
$ws.Range("L62").Value = 5425.75
$ws.Range("M62").Value = -4945.4443
$ws.Range("N62").Value = -6673.75
$ws.Range("H64").Value = 560318.9399999999
$ws.Range("I64").Value = 772783.1
$ws.Range("J64").Value = 7912
$ws.Range("K64").Value = 772783.1
$ws.Range("L64").Value = 7912
$ws.Range("M64").Value = -772535.1
$ws.Range("N64").Value = -8408
$ws.Range("H65").Value = 5501.8237
$ws.Range("I65").Value = 5569.4443
$ws.Range("J65").Value = 5425.75
$ws.Range("K65").Value = 27847.2215
$ws.Range("L65").Value = 27128.75
$ws.Range("M65").Value = -24727.2215
$ws.Range("N65").Value = -33368.75
$ws.Range("H67").Value = 560318.9399999999
$ws.Range("I67").Value = 772783.1
$ws.Range("J67").Value = 7912
$ws.Range("K67").Value = 772783.1
$ws.Range("L67").Value = 7912
$ws.Range("M67").Value = -771925.1
$ws.Range("N67").Value = -9628

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 11468.667
$ws.Range("I41").Value = 1762.4
$ws.Range("K41").Value = 1762.4
$ws.Range("M41").Value = -1348.4
$ws.Range("H61").Value = 3714.0435
$ws.Range("I61").Value = 2339.6365
$ws.Range("J61").Value = 4973.9165
$ws.Range("K61").Value = 2339.6365
$ws.Range("L61").Value = 4973.9165
$ws.Range("M61").Value = -2127.6365
$ws.Range("N61").Value = -5397.9165
$ws.Range("H136").Value = 3714.0435
$ws.Range("I136").Value = 2339.6365
$ws.Range("J136").Value = 4973.9165
$ws.Range("K136").Value = 7018.9095
$ws.Range("L136").Value = 14921.7495
$ws.Range("M136").Value = -4468.9095
$ws.Range("N136").Value = -20021.7495

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3288.611
$ws.Range("I105").Value = 3366.25
$ws.Range("K105").Value = 3366.25
$ws.Range("M105").Value = -1619.25
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H130").Value = 41250
$ws.Range("J130").Value = 41250
$ws.Range("L130").Value = 41250
$ws.Range("N130").Value = -51290
$ws.Range("H134").Value = 4550.2
$ws.Range("I134").Value = 3071.3635
$ws.Range("J134").Value = 6357.6665
$ws.Range("K134").Value = 9214.0905
$ws.Range("L134").Value = 19072.9995
$ws.Range("M134").Value = -6679.0905
$ws.Range("N134").Value = -24142.9995

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H99").Value = 12501582
$ws.Range("I99").Value = 12501582
$ws.Range("K99").Value = 12501582
$ws.Range("M99").Value = -12500084
$ws.Range("H105").Value = 508.07144
$ws.Range("I105").Value = 551
$ws.Range("K105").Value = 551
$ws.Range("M105").Value = 1196
$ws.Range("H107").Value = 396.06668
$ws.Range("I107").Value = 190.09091
$ws.Range("J107").Value = 962.5
$ws.Range("K107").Value = 190.09091
$ws.Range("L107").Value = 962.5
$ws.Range("M107").Value = 1729.90909
$ws.Range("N107").Value = -4802.5
$ws.Range("H126").Value = 12501582
$ws.Range("I126").Value = 12501582
$ws.Range("K126").Value = 37504746
$ws.Range("M126").Value = -37502276

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 323.25
$ws.Range("I8").Value = 323.25
$ws.Range("K8").Value = 969.75
$ws.Range("M8").Value = -830.75
$ws.Range("H12").Value = 27.947369
$ws.Range("I12").Value = 10.090909
$ws.Range("J12").Value = 52.5
$ws.Range("K12").Value = 30.272727
$ws.Range("L12").Value = 157.5
$ws.Range("M12").Value = 142.727273
$ws.Range("N12").Value = -503.5
$ws.Range("H107").Value = 360.88235
$ws.Range("I107").Value = 319.75
$ws.Range("J107").Value = 459.6
$ws.Range("K107").Value = 959.25
$ws.Range("L107").Value = 1378.8
$ws.Range("M107").Value = 960.75
$ws.Range("N107").Value = -5218.8
$ws.Range("H131").Value = 1657.279
$ws.Range("I131").Value = 664.6
$ws.Range("J131").Value = 1787.8948
$ws.Range("K131").Value = 1993.8
$ws.Range("L131").Value = 5363.6844
$ws.Range("M131").Value = 3046.2
$ws.Range("N131").Value = -15443.6844

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2779134.8
$ws.Range("I122").Value = 11111111
$ws.Range("J122").Value = 1809.3334
$ws.Range("K122").Value = 33333333
$ws.Range("L122").Value = 5428.0002
$ws.Range("M122").Value = -33330883
$ws.Range("N122").Value = -10328.0002
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3714.1082
$ws.Range("I136").Value = 1702.8572
$ws.Range("J136").Value = 9971.333000000001
$ws.Range("K136").Value = 5108.571599999999
$ws.Range("L136").Value = 29913.999
$ws.Range("M136").Value = -2558.571599999999
$ws.Range("N136").Value = -35013.999
$ws.Range("H139").Value = 53666.668
$ws.Range("J139").Value = 53666.668
$ws.Range("L139").Value = 53666.668
$ws.Range("N139").Value = -63946.668
$ws.Range("H141").Value = 42500
$ws.Range("J141").Value = 42500
$ws.Range("L141").Value = 42500
$ws.Range("N141").Value = -52860

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H131").Value = 62467.727
$ws.Range("J131").Value = 62467.727
$ws.Range("L131").Value = 62467.727
$ws.Range("N131").Value = -72547.727
$ws.Range("H136").Value = 9554055
$ws.Range("I136").Value = 15921462
$ws.Range("J136").Value = 2944.4285
$ws.Range("K136").Value = 47764386
$ws.Range("L136").Value = 8833.2855
$ws.Range("M136").Value = -47761836
$ws.Range("N136").Value = -13933.2855
